$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (legacy password hash); unprotect so the cell
# writes below (on an otherwise locked sheet) are allowed.
$ws.Unprotect()

# Update the confidential disclaimer date text (shared string used in A16)
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03064753504348171
$ws.Range("E2").Value = -0.004324324324324391

$ws.Range("D3").Value = 0.02404077155223762
$ws.Range("E3").Value = 0.008976225133430438

$ws.Range("D4").Value = 0.05283177143310964
$ws.Range("E4").Value = 0.009103977000479047

$ws.Range("D5").Value = 0.1381635180478894
$ws.Range("E5").Value = -0.001188253267696449

$ws.Range("D6").Value = 0.03162860362149068
$ws.Range("E6").Value = -0.01467351430667652

$ws.Range("D7").Value = 0.1191901530794584
$ws.Range("E7").Value = -0.00912336374454592

$ws.Range("D8").Value = 0.1019733054044395
$ws.Range("E8").Value = -0.004806767929244393

$ws.Range("D9").Value = 0.02820057181431918
$ws.Range("E9").Value = -0.004712746858168804

$ws.Range("D10").Value = 0.1227514593177937
$ws.Range("E10").Value = -0.01112009704811967

$ws.Range("D11").Value = 0.2490246396378509
$ws.Range("E11").Value = -0.007006248816511906

$ws.Range("D12").Value = 0.1015476710479292
$ws.Range("E12").Value = 0.004777253130774151

$ws.Range("E13").Value = -0.004399128380256467

$wb.Save()
